$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("1E - Constant")

# Move the active selection from H30 to H20.
$ws.Activate()
$ws.Range("H20").Select()

# H20 becomes a formula instead of a hard-coded value.
$ws.Range("H20").Formula = "=3076.92/60"

# I20 gains a formula; I21:I29 gain plain numeric values (second trial run).
$ws.Range("I20").Formula = "=3099.55/60"
$ws.Range("I21").Value = 7249.76
$ws.Range("I22").Value = 6812.75
$ws.Range("I23").Value = 6209.46
$ws.Range("I24").Value = 5994.1
$ws.Range("I25").Value = 5885.34
$ws.Range("I26").Value = 5815.4
$ws.Range("I27").Value = 5774.0110000000004
$ws.Range("I28").Value = 5720.67
$ws.Range("I29").Value = 5682.38
